$d = $word.ActiveDocument

# Locate the end of the "LOQ4053..." requirement paragraph text - the
# three trailing paragraphs (the blank spacer, the "Ver no Jupiter..."
# line, and the "(c) 2020 ..." footer line) must be removed, leaving the
# requirement paragraph followed directly by the pre-existing blank
# paragraph / page break.
$anchorStart = $d.Content
$anchorStart.Find.Execute("LOQ4053: Balanços de Massa e Energia (Requisito fraco)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$afterRequisito = $anchorStart.End

# Locate the end of the copyright/footer paragraph text.
$anchorEnd = $d.Content
$anchorEnd.Find.Execute("© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$afterCopyright = $anchorEnd.End

# Expand by one character on each side to swallow the paragraph marks:
# one to skip past the end-of-paragraph mark that closes the requirement
# line, and one to include the end-of-paragraph mark that closes the
# copyright line - this removes the three intervening paragraphs
# (blank spacer, "Ver no Jupiter..." line, copyright line) as a block.
$deleteRange = $d.Range($afterRequisito + 1, $afterCopyright + 1)
$deleteRange.Delete()
